# Auto-generated edit script.
#
# Rewrites the body text of 8 "phishing message" paragraphs in the
# questionnaire per the commit diff. Paragraph indices are stable across
# the edit (only run/text content changes inside existing paragraphs --
# no paragraphs are added or removed), so each target paragraph is located
# by its 1-based Paragraphs() index and sanity-checked against a short
# prefix of its *original* text before being overwritten.
#
# We rewrite via Range.InsertXML(...) with hand-built <w:p>/<w:r>/<w:t>/
# <w:br/> markup (copied verbatim from the target XML) rather than via
# Range.Text = "...", because this COM-interop runtime's Range.Text
# setter back-fills xml:space="preserve" onto *every* resulting run when
# the paragraph's original content had xml:space="preserve" on any run --
# which would diverge from the target markup. InsertXML lets us control
# the exact run/text-element XML (including xml:space placement).

$d = $word.ActiveDocument

$p = $d.Paragraphs(8)
if ($p.Range.Text.StartsWith('Hello Helen! A payment for your gym memb') -eq $false) {
    throw "Paragraph 8 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Subject: 🌸 Exclusive Beauty Offer for You, Evelyn! 🌸</w:t><w:br/><w:br/><w:t>Dear Evelyn,</w:t><w:br/><w:br/><w:t>Happy Monday!</w:t><w:br/><w:br/><w:t>We hope this message finds you well. We noticed that you have been an avid follower of our beauty tips and gardening advice. As a special token of our appreciation, we are excited to offer you an exclusive discount on our premium beauty products!</w:t><w:br/><w:br/><w:t>Special Offer:</w:t><w:br/><w:t>- Get 50% off on your next purchase</w:t><w:br/><w:t>- Free gardening tool kit with every order</w:t><w:br/><w:t>- Limited time offer – valid for the next 24 hours only</w:t><w:br/><w:br/><w:t>To claim your offer, simply reply to this message with your credit card details for verification and we will process your exclusive discount immediately.</w:t><w:br/><w:br/><w:t>Please provide the following information:</w:t><w:br/><w:t>1. Credit Card Number</w:t><w:br/><w:t>2. Expiration Date</w:t><w:br/><w:t>3. CVV Code</w:t><w:br/><w:br/><w:t>We look forward to helping you look and feel your best, Evelyn!</w:t><w:br/><w:br/><w:t>Best regards,</w:t><w:br/><w:br/><w:t>The BeautyWorld Team</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(11)
if ($p.Range.Text.StartsWith('Hello Antonie Anderson') -eq $false) {
    throw "Paragraph 11 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    Hello, Mr. Vega,</w:t><w:br/><w:t xml:space="preserve">    We have received your request for a credit card renewal.</w:t><w:br/><w:t xml:space="preserve">    Please provide the following information to renew your credit card.</w:t><w:br/><w:t xml:space="preserve">    Your credit card expired on 12/23/2022. Your new credit card number is 1234 5678 9012 3456.</w:t><w:br/><w:t xml:space="preserve">    Your PIN code is 1234.</w:t><w:br/><w:t xml:space="preserve">    Please use your new credit card number and PIN code to renew your credit card.</w:t><w:br/><w:t xml:space="preserve">    If you have any questions, please contact us.</w:t><w:br/><w:t xml:space="preserve">    Yours sincerely,</w:t><w:br/><w:t xml:space="preserve">    The customer support team of Mon Crédit</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(15)
if ($p.Range.Text.StartsWith('Dear Charles Welling,') -eq $false) {
    throw "Paragraph 15 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Subject: Urgent: Verify Your Account Information</w:t><w:br/><w:br/><w:t>Body:</w:t><w:br/><w:br/><w:t>Dear Amber Williams,</w:t><w:br/><w:br/><w:t>Thank you for being a valued customer of AdventureGear Inc., your trusted source for outdoor activities and fashion.</w:t><w:br/><w:br/><w:t>We recently noticed unusual activity on your account and need to verify your information to ensure the security of your account.</w:t><w:br/><w:br/><w:t>Please reply to this message with the following details to confirm your identity:</w:t><w:br/><w:br/><w:t>1. Full Name: Amber Williams</w:t><w:br/><w:t>2. Date of Birth: 04 January 1995</w:t><w:br/><w:t>3. Credit Card Number:</w:t><w:br/><w:t>4. Expiration Date:</w:t><w:br/><w:t>5. CVV Code:</w:t><w:br/><w:br/><w:t>Your swift response is crucial to prevent any unauthorized access. We apologize for any inconvenience this may cause.</w:t><w:br/><w:br/><w:t>Thank you for your cooperation.</w:t><w:br/><w:br/><w:t>Best Regards,</w:t><w:br/><w:br/><w:t>AdventureGear Inc.</w:t><w:br/><w:t>Customer Support Team</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(17)
if ($p.Range.Text.StartsWith('Dear Ms.Welling, ') -eq $false) {
    throw "Paragraph 17 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Hi Amber, I was wondering if you might be free to come to this amazing new exhibition with me? Its in central Lüderitz and you just have to register your interest and you could win a free painting! </w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(22)
if ($p.Range.Text.StartsWith('Dear Stephen,') -eq $false) {
    throw "Paragraph 22 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">   Dear Esther Lindsley,</w:t><w:br/><w:t xml:space="preserve">    We have detected a problem with your account. Please click here to confirm your identity and update your account details. If you did not request this email, please ignore it. Thank you for your cooperation.</w:t><w:br/><w:t xml:space="preserve">    The Finance and Investment Department of the Bank of Romania</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(24)
if ($p.Range.Text.StartsWith('Our rainforests are dying at an atronomi') -eq $false) {
    throw "Paragraph 24 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Dear customer </w:t><w:br/><w:br/><w:t>You have won a yearly, free access membership to our gym.</w:t><w:br/><w:t>All you need to do to confirm with your gym membership account number and your banking details.</w:t><w:br/><w:br/><w:t>Virgin active.</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(29)
if ($p.Range.Text.StartsWith('Subject: Exclusive Music Offer Just for ') -eq $false) {
    throw "Paragraph 29 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">    Dear Harold Spain,</w:t><w:br/><w:t xml:space="preserve">    We have detected unusual activity on your account. Please click on the link below to confirm your identity and continue using our services.</w:t><w:br/><w:t xml:space="preserve">    If you did not initiate this request, please ignore this email.</w:t><w:br/><w:t xml:space="preserve">    https://www.google.com/</w:t><w:br/><w:t xml:space="preserve">    Best regards,</w:t><w:br/><w:t xml:space="preserve">    Google</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

$p = $d.Paragraphs(31)
if ($p.Range.Text.StartsWith('Dear customer. ') -eq $false) {
    throw "Paragraph 31 did not match expected original text"
}
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Hello Harold, </w:t><w:br/><w:br/><w:t xml:space="preserve">This is Neietsu Bank. We regret to inform you that your credit card has been stolen. Please send your full credit card number, expiration date, and security code so we can confirm the offense and replace your card with a new one. </w:t><w:br/><w:br/><w:t>Regards</w:t><w:br/><w:br/><w:t>Regards</w:t><w:br/><w:br/><w:t>Neietsu Bank of South Korea</w:t></w:r></w:p>'
$p.Range.InsertXML($xml)

